# "Moved logics to 3E to fix corruption bug"
#
# The logicEntry / logicFile allocation rows that used to live at the end of
# BANK0x3C are moved to BANK0x3D's sibling bank, BANK0x3E (which only had a
# single "logicEntryAddresses" row before). BANK0x3E's first row is renamed
# to "pictureFile" and enlarged, and the two moved rows are appended after it
# with fresh (non-shared) formulas that chain from BANK0x3E's own rows
# instead of BANK0x3C's.

$wb = $excel.ActiveWorkbook

$ws3e = $wb.Worksheets.Item("BANK0x3E")
$ws3c = $wb.Worksheets.Item("BANK0x3C")
$wsDyn = $wb.Worksheets.Item(1)

# --- BANK0x3E: rename/resize the existing allocation row, then add the two
#     logic rows that used to sit at the bottom of BANK0x3C ---
$ws3e.Range("A2").Value = "pictureFile"
$ws3e.Range("C2").Value = 6

$ws3e.Range("A3").Value = "logicEntry"
$ws3e.Range("B3").Formula = "=E2 + 1"
$ws3e.Range("C3").Value = 8
$ws3e.Range("D3").Value = 255
$ws3e.Range("E3").Formula = "=C3*D3"

$ws3e.Range("A4").Value = "logicFile"
$ws3e.Range("B4").Formula = "=B3+ E3+1"
$ws3e.Range("C4").Value = 9
$ws3e.Range("D4").Value = 255
$ws3e.Range("E4").Formula = "=C4*D4"

# --- BANK0x3C: the logicEntry/logicFile rows (7 & 8) are gone now, so clear
#     them, and repoint the next row's running-total start at the rows that
#     replaced them on BANK0x3E ---
$ws3c.Range("A7:E8").ClearContents()
$ws3c.Range("B9").Formula = "=BANK0x3E!B4+ BANK0x3E!E4+1"

# --- view/selection bookkeeping, matching what the author last had on
#     screen: the dynamic-bank sheet's selection moved, and the active tab
#     moved from BANK0x3C to BANK0x3E ---
$wsDyn.Activate()
$wsDyn.Range("E6").Select()

$ws3c.Activate()
$ws3c.Range("A7:E8").Select()

$ws3e.Activate()
$ws3e.Range("B3").Select()
